$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Report-" test-scenario placeholder text is replaced by a full,
# specific filename-prefix scenario used for the non-oncology import tool.
$ws.Range("J4").Value = "StandardExcelReport-NewImportLogic_3 - Test_Automation_3-Clinical-2023_"

# The active selection moves on to the newly edited cell, scrolling the
# view one column to the right (H1 -> I1 becomes the top-left visible cell).
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("J4").Select()

# Column J is widened to accommodate the longer replacement text.
$ws.Columns.Item(10).ColumnWidth = 52.333333333333336
